$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the view: scroll so row 4 is the top-left visible row, and move
# the active selection to J6.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("J6").Select()

# Clear the "1" markers at several cells (set them to 0) - these appear to
# represent animation/collision points removed since animations now live
# in UCS.
$ws.Range("T5").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("Q20").Value = 0
